$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75 (pushes existing rows 75-179 down to 76-180)
$ws.Rows("75").Insert()

# Populate the new row 75 with the new record's data.
# Columns A,B,C,E,F,G,N,Q,R are constant across this block of rows.
$ws.Range("A75").Value = 5
$ws.Range("B75").Value = "Macroferia Regional de Talca"
$ws.Range("C75").Value = "Maule"
$ws.Range("D75").Value2 = 44580
$ws.Range("E75").Value = 7
$ws.Range("F75").Value = 100112045
$ws.Range("G75").Value = "Zapallo"
$ws.Range("H75").Value = "Camote"
$ws.Range("I75").Value = "1a nueva(o)"
$ws.Range("J75").Value = 600
$ws.Range("K75").Value = 300
$ws.Range("L75").Value = 400
$ws.Range("M75").Value = 367
$ws.Range("N75").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O75").Value = "Región del Maule"
$ws.Range("P75").Value = 367
$ws.Range("Q75").Value = 1
$ws.Range("R75").Value = "Hortaliza"
